$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like text label to a cell without letting Excel's
# automatic date recognition turn it into a date serial number. We build
# the text via a temporary formula cell, copy it, and paste-special just
# the resulting value (which keeps it as plain text) into the destination.
function Set-TextLabel($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# --- Update existing row 193 with revised figures ---
$ws.Range("E193").Value = -0.22
$ws.Range("H193").Value = 2.07
$ws.Range("I193").Value = 1.99
$ws.Range("L193").Value = 6.86
$ws.Range("M193").Value = 11.08
$ws.Range("N193").Value = 7.72
$ws.Range("O193").Value = 7.24
$ws.Range("P193").Value = 6.4

# --- Append new rows 194-198 with new daily data ---

# Row 194: 28-09-2021
$ws.Range("A194").Value = "28-09-2021"
$ws.Range("B194").Value = 1.54
$ws.Range("C194").Value = 0.99
$ws.Range("D194").Value = 0.09
$ws.Range("E194").Value = -0.2
$ws.Range("G194").Value = 3.45
$ws.Range("H194").Value = 2.12
$ws.Range("I194").Value = 1.99
$ws.Range("J194").Value = 7.31
$ws.Range("K194").Value = 1.83
$ws.Range("L194").Value = 6.93
$ws.Range("M194").Value = 11.12
$ws.Range("N194").Value = 7.71
$ws.Range("O194").Value = 7.33
$ws.Range("P194").Value = 6.42

# Row 195: 29-09-2021
$ws.Range("A195").Value = "29-09-2021"
$ws.Range("B195").Value = 1.52
$ws.Range("C195").Value = 0.99
$ws.Range("D195").Value = 0.08
$ws.Range("E195").Value = -0.21
$ws.Range("F195").Value = 2.27
$ws.Range("G195").Value = 3.42
$ws.Range("H195").Value = 2.15
$ws.Range("I195").Value = 1.98
$ws.Range("J195").Value = 7.29
$ws.Range("K195").Value = 1.83
$ws.Range("L195").Value = 6.79
$ws.Range("M195").Value = 11.07
$ws.Range("N195").Value = 7.68
$ws.Range("O195").Value = 7.36
$ws.Range("P195").Value = 6.43

# Row 196: 30-09-2021
$ws.Range("A196").Value = "30-09-2021"
$ws.Range("B196").Value = 1.49
$ws.Range("C196").Value = 1.02
$ws.Range("D196").Value = 0.07000000000000001
$ws.Range("E196").Value = -0.2
$ws.Range("F196").Value = 2.25
$ws.Range("G196").Value = 3.38
$ws.Range("H196").Value = 2.16
$ws.Range("I196").Value = 2.06
$ws.Range("J196").Value = 7.32
$ws.Range("K196").Value = 1.86
$ws.Range("L196").Value = 6.79
$ws.Range("M196").Value = 11.11
$ws.Range("N196").Value = 7.66
$ws.Range("O196").Value = 7.37
$ws.Range("P196").Value = 6.44

# Row 197: 01-10-2021 (would otherwise be auto-parsed as a date, so use
# the text-preserving helper)
Set-TextLabel $ws.Range("A197") "01-10-2021"
$ws.Range("B197").Value = 1.46
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 0.06
$ws.Range("E197").Value = -0.22
$ws.Range("G197").Value = 3.38
$ws.Range("H197").Value = 2.21
$ws.Range("I197").Value = 2.08
$ws.Range("J197").Value = 7.35
$ws.Range("K197").Value = 1.89
$ws.Range("L197").Value = 6.82
$ws.Range("M197").Value = 10.99
$ws.Range("N197").Value = 7.64
$ws.Range("O197").Value = 7.4
$ws.Range("P197").Value = 6.39

# Row 198: 04-10-2021 (same date-parsing concern as row 197)
Set-TextLabel $ws.Range("A198") "04-10-2021"
$ws.Range("B198").Value = 1.5
$ws.Range("C198").Value = 1.02
$ws.Range("D198").Value = 0.05
$ws.Range("E198").Value = -0.21
$ws.Range("G198").Value = 3.38
$ws.Range("H198").Value = 2.26
$ws.Range("I198").Value = 2.11
$ws.Range("J198").Value = 7.42
$ws.Range("K198").Value = 1.88
$ws.Range("L198").Value = 6.74
$ws.Range("O198").Value = 7.39
$ws.Range("P198").Value = 6.37
